$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '26.553.34'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.11%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '1.811.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -0.16%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'" + '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.48%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -0.36%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '306.00'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.86%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '0.4549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.40%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.3595'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -1.99%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + '46.33'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  +2.19%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.07109'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -0.38%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.8924'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +1.28%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '0.07711'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -0.47%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '19.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -0.22%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '1.816.83'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +0.18%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '5.255'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -0.87%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '6.293'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -1.28%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + '86.64'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -0.19%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D19').Value = "'" + '0.000008556'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -0.46%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  -0.38%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '26.581.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -0.01%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -0.89%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '4.961'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -1.19%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '10.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.18%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '1.927'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -2.77%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '151.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +0.29%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + '17.81'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -0.78%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '2.017'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -3.04%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '112.46'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -0.61%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '4.832'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  -0.63%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '0.08722'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +0.31%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '3.125'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +3.04%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '0.7388'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +0.95%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'" + 'Filecoin'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'" + '4.435'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -1.49%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'" + 'RenderToken'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'" + '2.721'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +2.24%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '1.110'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -0.86%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '1.071'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -1.33%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -1.39%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + '2.913'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +0.72%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'" + '0.05071'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -1.17%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.5081'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  +1.47%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '6.779'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -2.99%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -3.23%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '8.008'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -1.84%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '0.4685'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +1.76%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '1.003'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -0.34%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '9.939'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  -0.30%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '99.24'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -2.13%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -1.41%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -0.07%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -1.47%  '
$ws.Range('E51').Style = 'Normal'
